$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..48 down to 4..49
$ws.Rows("3:3").Insert()

# Populate the new row 3 with values (same pattern as the surrounding rows,
# with new Fecha / Volumen / Precio minimo / Precio maximo / Precio promedio
# ponderado / Precio $/Kg)
$ws.Range("A3").Value = 6
$ws.Range("B3").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44812
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112035
$ws.Range("G3").Value = "Bruselas (repollito)"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 410
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17488
$ws.Range("N3").Value = "$/malla 15 kilos"
$ws.Range("O3").Value = "Provincia de Quillota"
$ws.Range("P3").Value = 1166
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = "Hortaliza"
